# "changes in concise marksheet" - update Corr/total marks
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Marking row: marks awarded per correct answer (B11)
$ws.Range("B11").Value = 5

# Total row: correct marks obtained (B12) and "obtained/max" summary (E12)
$ws.Range("B12").Value = 100
$ws.Range("E12").Value = "100/140"
